# Update the "想去人数" (number of people interested) column (F) for
# several rows on the "展览" sheet and the "全部类型" sheet, reflecting a
# fresh scrape of the source data (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3542
$ws1.Range("F5").Value = 3542
$ws1.Range("F8").Value = 513
$ws1.Range("F14").Value = 30
$ws1.Range("F15").Value = 691
$ws1.Range("F16").Value = 312
$ws1.Range("F22").Value = 4893
$ws1.Range("F32").Value = 4442
$ws1.Range("F36").Value = 1001
$ws1.Range("F40").Value = 864
$ws1.Range("F41").Value = 968
$ws1.Range("F42").Value = 1230

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3542
$ws4.Range("F8").Value = 3542
$ws4.Range("F11").Value = 513
$ws4.Range("F17").Value = 30
$ws4.Range("F18").Value = 691
$ws4.Range("F19").Value = 312
$ws4.Range("F26").Value = 4893
$ws4.Range("F36").Value = 4442
$ws4.Range("F41").Value = 1001
$ws4.Range("F45").Value = 864
$ws4.Range("F46").Value = 968
$ws4.Range("F48").Value = 1232
